$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("G2").Value = 7.132470666666666
    $ws.Range("H2").Value = 21.397412
    $ws.Range("I2").Value = 0.1078130252899183
    $ws.Range("J2").Value = 0.1078130252899183
    $ws.Range("M2").Value = 31.61061466666667
    $ws.Range("N2").Value = 94.831844
    $ws.Range("O2").Value = 0.8860472269592234
    $ws.Range("P2").Value = 0.8860472269592234
    $ws.Range("Q2").Value = 225.4617818653031
    $ws.Range("R2").Value = 2029.156036787728
    $ws.Range("S2").Value = 0.09552743208821676
    $ws.Range("T2").Value = 0.09552743208821676
    $ws.Range("G3").Value = 7.132470666666666
    $ws.Range("H3").Value = 21.397412
    $ws.Range("I3").Value = 0.1078130252899183
    $ws.Range("J3").Value = 0.1078130252899183
    $ws.Range("O3").Value = 0.04688826274109129
    $ws.Range("P3").Value = 0.04688826274109129
    $ws.Range("Q3").Value = 11.93109232163022
    $ws.Range("R3").Value = 107.379830894672
    $ws.Range("S3").Value = 0.005055165456705611
    $ws.Range("T3").Value = 0.00505516545670561
    $ws.Range("G4").Value = 7.132470666666666
    $ws.Range("H4").Value = 21.397412
    $ws.Range("I4").Value = 0.1078130252899183
    $ws.Range("J4").Value = 0.1078130252899183
    $ws.Range("M4").Value = 2.392593
    $ws.Range("N4").Value = 7.177778999999999
    $ws.Range("O4").Value = 0.06706451029968528
    $ws.Range("P4").Value = 0.06706451029968527
    $ws.Range("Q4").Value = 17.065099389772
    $ws.Range("R4").Value = 153.585894507948
    $ws.Range("S4").Value = 0.007230427744995958
    $ws.Range("T4").Value = 0.007230427744995957
    $ws.Range("I5").Value = 0.2490596131114117
    $ws.Range("J5").Value = 0.2490596131114118
    $ws.Range("M5").Value = 31.61061466666667
    $ws.Range("N5").Value = 94.831844
    $ws.Range("O5").Value = 0.8860472269592234
    $ws.Range("P5").Value = 0.8860472269592234
    $ws.Range("Q5").Value = 520.8408168844219
    $ws.Range("R5").Value = 4687.567351959797
    $ws.Range("S5").Value = 0.2206785795449034
    $ws.Range("T5").Value = 0.2206785795449034
    $ws.Range("I6").Value = 0.2490596131114117
    $ws.Range("J6").Value = 0.2490596131114118
    $ws.Range("O6").Value = 0.04688826274109129
    $ws.Range("P6").Value = 0.04688826274109129
    $ws.Range("S6").Value = 0.01167797257776242
    $ws.Range("T6").Value = 0.01167797257776242
    $ws.Range("I7").Value = 0.2490596131114117
    $ws.Range("J7").Value = 0.2490596131114118
    $ws.Range("M7").Value = 2.392593
    $ws.Range("N7").Value = 7.177778999999999
    $ws.Range("O7").Value = 0.06706451029968528
    $ws.Range("P7").Value = 0.06706451029968527
    $ws.Range("Q7").Value = 39.422203767079
    $ws.Range("R7").Value = 354.799833903711
    $ws.Range("S7").Value = 0.0167030609887459
    $ws.Range("T7").Value = 0.0167030609887459
    $ws.Range("G8").Value = 42.546687
    $ws.Range("H8").Value = 127.640061
    $ws.Range("I8").Value = 0.6431273615986699
    $ws.Range("J8").Value = 0.6431273615986699
    $ws.Range("M8").Value = 31.61061466666667
    $ws.Range("N8").Value = 94.831844
    $ws.Range("O8").Value = 0.8860472269592234
    $ws.Range("P8").Value = 0.8860472269592234
    $ws.Range("Q8").Value = 1344.926928100276
    $ws.Range("R8").Value = 12104.34235290248
    $ws.Range("S8").Value = 0.5698412153261032
    $ws.Range("T8").Value = 0.5698412153261032
    $ws.Range("G9").Value = 42.546687
    $ws.Range("H9").Value = 127.640061
    $ws.Range("I9").Value = 0.6431273615986699
    $ws.Range("J9").Value = 0.6431273615986699
    $ws.Range("O9").Value = 0.04688826274109129
    $ws.Range("P9").Value = 0.04688826274109129
    $ws.Range("Q9").Value = 71.17147399552401
    $ws.Range("R9").Value = 640.543265959716
    $ws.Range("S9").Value = 0.03015512470662326
    $ws.Range("T9").Value = 0.03015512470662325
    $ws.Range("G10").Value = 42.546687
    $ws.Range("H10").Value = 127.640061
    $ws.Range("I10").Value = 0.6431273615986699
    $ws.Range("J10").Value = 0.6431273615986699
    $ws.Range("M10").Value = 2.392593
    $ws.Range("N10").Value = 7.177778999999999
    $ws.Range("O10").Value = 0.06706451029968528
    $ws.Range("P10").Value = 0.06706451029968527
    $ws.Range("Q10").Value = 101.796905489391
    $ws.Range("R10").Value = 916.1721494045188
    $ws.Range("S10").Value = 0.04313102156594342
    $ws.Range("T10").Value = 0.04313102156594341
